$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextValue 'D2' '50.979.96'
$ws.Range('E2').Value = '  -0.95%  '
Set-TextValue 'D3' '2.933.11'
$ws.Range('E3').Value = '  -1.49%  '
$ws.Range('E4').Value = '  +0.05%  '
Set-TextValue 'D5' '374.28'
$ws.Range('E5').Value = '  -1.63%  '
Set-TextValue 'D6' '101.58'
$ws.Range('E6').Value = '  -2.96%  '
Set-TextValue 'D7' '0.535'
$ws.Range('E7').Value = '  -1.73%  '
$ws.Range('E8').Value = '  +0.06%  '
Set-TextValue 'D9' '0.581'
$ws.Range('E9').Value = '  -2.58%  '
Set-TextValue 'D10' '36.36'
$ws.Range('E10').Value = '  -2.32%  '
$ws.Range('E11').Value = '  -0.58%  '
Set-TextValue 'D12' '0.0833'
$ws.Range('E12').Value = '  -1.31%  '
Set-TextValue 'D13' '3.393.10'
$ws.Range('E13').Value = '  -1.57%  '
Set-TextValue 'D14' '17.89'
$ws.Range('E14').Value = '  -3.54%  '
Set-TextValue 'D15' '7.32'
$ws.Range('E15').Value = '  -2.25%  '
Set-TextValue 'D16' '2.928.93'
$ws.Range('E16').Value = '  -1.49%  '
Set-TextValue 'D17' '0.972'
$ws.Range('E17').Value = '  -0.33%  '
Set-TextValue 'D18' '50.929.07'
$ws.Range('E18').Value = '  -1.01%  '
Set-TextValue 'D19' '3.14'
$ws.Range('E19').Value = '  -6.67%  '
Set-TextValue 'D20' '7.13'
$ws.Range('E20').Value = '  -3.54%  '
Set-TextValue 'D21' '12.50'
$ws.Range('E21').Value = '  -3.69%  '
Set-TextValue 'D22' '0.0₃0953'
$ws.Range('E22').Value = '  -0.98%  '
Set-TextValue 'D23' '263.67'
$ws.Range('E23').Value = '  +0.34%  '
Set-TextValue 'D24' '68.13'
$ws.Range('E24').Value = '  -1.39%  '
Set-TextValue 'D25' '2.89'
$ws.Range('E25').Value = '  +2.61%  '
Set-TextValue 'D26' '7.74'
Set-TextValue 'D27' '8.02'
$ws.Range('E27').Value = '  +7.60%  '
$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D28' '1.00'
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D29' '0.167'
$ws.Range('E29').Value = '  -1.77%  '
$ws.Range('E30').Value = '  -0.82%  '
Set-TextValue 'D31' '25.60'
$ws.Range('E31').Value = '  -1.60%  '
Set-TextValue 'D32' '9.82'
$ws.Range('E32').Value = '  -0.28%  '
$ws.Range('E33').Value = '  -0.61%  '
Set-TextValue 'D34' '0.0452'
$ws.Range('E34').Value = '  -0.68%  '
Set-TextValue 'D35' '33.44'
$ws.Range('E35').Value = '  -4.17%  '
Set-TextValue 'D36' '2.01'
$ws.Range('E36').Value = '  -3.46%  '
$ws.Range('E37').Value = '  -0.13%  '
Set-TextValue 'D38' '2.96'
$ws.Range('E38').Value = '  -4.29%  '
Set-TextValue 'D39' '2.53'
$ws.Range('E39').Value = '  -2.48%  '
$ws.Range('E40').Value = '  -1.73%  '
Set-TextValue 'D41' '16.29'
$ws.Range('E41').Value = '  -6.63%  '
Set-TextValue 'D42' '1.78'
$ws.Range('E42').Value = '  -4.07%  '
Set-TextValue 'D43' '121.62'
$ws.Range('E43').Value = '  -1.84%  '
Set-TextValue 'D44' '20.82'
$ws.Range('E44').Value = '  -6.28%  '
$ws.Range('E45').Value = '  -1.75%  '
Set-TextValue 'D46' '0.272'
$ws.Range('E46').Value = '  -3.09%  '
Set-TextValue 'D47' '2.30'
$ws.Range('E47').Value = '  -4.48%  '
Set-TextValue 'D48' '3.20'
$ws.Range('E48').Value = '  -1.09%  '
Set-TextValue 'D49' '1.974.48'
$ws.Range('E49').Value = '  -3.11%  '
Set-TextValue 'D50' '0.0344'
$ws.Range('E50').Value = '  -3.33%  '
Set-TextValue 'D51' '5.02'
